# Apply the changes described by the commit "Fix Errors and update PR template"
#
# 1. Rename "Sheet1" to "Purchase Request"
# 2. Update the Department field (I7) from "  IT-Bacolod" to "ITBS"
# 3. Update the Site PR field (I8) from "PR100112341" to "PR10011234112"
# 4. Leave the final selection on cell F20 (matches the saved cursor position)

$wb = $excel.ActiveWorkbook

# 1. Rename the main worksheet
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Purchase Request"

# 2 & 3. Correct the Department / Site PR values on the form
$ws.Range("I7").Value = "ITBS"
$ws.Range("I8").Value = "PR10011234112"

# 4. Leave the selection on F20, same as the saved workbook
$ws.Activate()
$ws.Range("F20").Select()

Write-Output "Done: renamed sheet, updated I7/I8, set selection to F20"
